# REQ-17 and REQ-18 are added to system test report
#
# Fills in the four previously-blank template rows (18-21) on the
# "Test Cases & Results" sheet with the new REQ-15 .. REQ-18 test cases,
# matching columns:
#   D = Requirement_ID, E = Priority, F = Description/Test Summary,
#   G = Pre-Condition, H = Test Steps, I = Expected Result,
#   J = Actual Result (mirrors Expected Result), K = Test Result (unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases & Results")
$ws.Activate()

# --- Row 18 (test case 16 / REQ-15) ---
$ws.Range("D18").Value = "REQ-15"
$ws.Range("E18").Value = "Low Impact"
$ws.Range("F18").Value = "Test after the servo motor returns to its original position, buzzer is beeped quickly for 3 seconds"
$ws.Range("G18").Value = "Servo motor must be back to its original position"
$ws.Range("H18").Value = "Follow the same steps as test case 15, no additional steps required"
$ws.Range("I18").Value = "Buzzer is beeped quickly for 3 seconds"
$ws.Range("J18").Value = "Buzzer is beeped quickly for 3 seconds"
$ws.Rows.Item(18).RowHeight = 57.6

# --- Row 19 (test case 17 / REQ-16) ---
$ws.Range("D19").Value = "REQ-16"
$ws.Range("E19").Value = "High Impact"
$ws.Range("F19").Value = 'Test that the LCD shows "Please Scan Your Card" after the REQ-15 is completed '
$ws.Range("G19").Value = "The buzzer has finished beeping and firebase is updated"
$ws.Range("H19").Value = 'Follow the same steps as test case 16, no additional steps required '
$ws.Range("I19").Value = 'LCD displays "Please Scan Your Card"'
$ws.Range("J19").Value = 'LCD displays "Please Scan Your Card"'
$ws.Rows.Item(19).RowHeight = 43.2

# --- Row 20 (test case 18 / REQ-17) ---
$ws.Range("D20").Value = "REQ-17"
$ws.Range("E20").Value = "High Impact"
$ws.Range("F20").Value = "Test that in REQ-04 if option 2 is selected on the matrix keypad, REQ-18 is started"
$ws.Rows.Item(20).RowHeight = 43.2

# --- Row 21 (test case 19 / REQ-18) ---
$ws.Range("D21").Value = "REQ-18"
$ws.Range("E21").Value = "High Impact"
$ws.Range("G21").Value = "The user selects 2 in the Matrix keypad from REQ-04"
$ws.Range("I21").Value = 'LCD displays "Scan book, 0 to end"'
$ws.Range("J21").Value = 'LCD displays "Scan book, 0 to end"'
$ws.Range("F21").Value = 'Test that the LCD shows "Scan book, 0 to end" after option 2 is selected'
$ws.Rows.Item(21).RowHeight = 43.2

# Reflect the author's final on-screen selection/scroll position.
$ws.Range("L20").Select()
